$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3.195530200612898
$ws.Range("E2").Value = 1.196910142261056
$ws.Range("F2").Value = 3.171945884630063
$ws.Range("G2").Value = 1.166983854326101
